# Adding a few more assets to the "Holdings" sheet (Index/Asset mapping table).
#
# Before: rows 2-18 held 17 asset/class pairs.
# After:  rows 2-26 hold 25 asset/class pairs -- 3 new "commodities" rows are
#         inserted right after the existing commodities block (old row 13),
#         and 5 new "Equities" rows are appended after the existing equities
#         block (old row 18 / new row 21).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Holdings")

# --- Insert 3 new "commodities" rows before the old row 13 (BIL / cash_equivalents) ---
$ws.Rows.Item(13).Insert()
$ws.Rows.Item(13).Insert()
$ws.Rows.Item(13).Insert()

$newCommodities = @(
    @("DBB", "commodities"),
    @("GSG", "commodities"),
    @("BNO", "commodities")
)
for ($i = 0; $i -lt $newCommodities.Length; $i++) {
    $r = 13 + $i
    $ws.Rows.Item($r).RowHeight = 18.75
    $ws.Range("A$r").Value = $newCommodities[$i][0]
    $ws.Range("B$r").Value = $newCommodities[$i][1]
}

# After the insertion above, the old rows 13-18 (BIL..QQQ) have shifted down
# to rows 16-21.

# --- Append 5 new "Equities" rows after the old QQQ row (now row 21) ---
$newEquities = @(
    @("EEM", "Equities"),
    @("DIA", "Equities"),
    @("IJR", "Equities"),
    @("IVE", "Equities"),
    @("ACWI", "Equities")
)
for ($i = 0; $i -lt $newEquities.Length; $i++) {
    $r = 22 + $i
    $ws.Rows.Item($r).RowHeight = 18.75
    $ws.Range("A$r").Value = $newEquities[$i][0]
    $ws.Range("B$r").Value = $newEquities[$i][1]
}

# The newly appended rows (22-26) fall outside the sheet's previous used
# range, so they don't automatically inherit the A/B column cell styles
# (s="3"/s="4"). Copy formatting from an existing Equities row (row 19,
# "VTI") onto them so they match the rest of the table.
$ws.Range("A19:B19").Copy()
$ws.Range("A22:B26").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
